# Automatische test-sync: 2025-07-29 22:04:50
# Adds the 17th test mail to the "Logs" sheet, updates the "Dashboard"
# pivot-style summary with the new "Planning / Afspraak" category, and
# extends the chart series range so the new row is plotted too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 19 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A19").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("B19").Value = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D19").Value = "Planning / Afspraak"
$logs.Range("E19").Value = "Beste [Naam],`nBedankt voor je e-mail. Ik heb de demo op vrijdag om 11:00 uur bij Van Dijk ingepland. Mocht er iets wijzigen of als er nog vragen zijn, laat het me gerust weten.`nMet vriendelijke groet,`n[Jouw Naam]"
$logs.Range("F19").Value = "2025-07-29 22:04:18"
$logs.Range("G19").Value = "Ja"
$logs.Range("H19").Value = "Nee"
$logs.Range("I19").Value = "Ja"
$logs.Range("J19").Value = "Nee"

# The multi-line "Antwoord" text would otherwise leave an explicit
# custom row height behind; auto-fit puts the row back to the sheet's
# default (non-custom) height, matching the other rows.
$logs.Rows.Item(19).AutoFit()

# Conditional formatting on these columns previously covered rows 2-18;
# extend each of them so row 19 is covered as well.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "18")
    $newRange = $logs.Range($col + "2:" + $col + "19")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: append the new category total (row 7)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Planning / Afspraak"
$dash.Range("B7").Value = 1

# ---------------------------------------------------------------------
# 3) Chart: extend category/value series references to row 7
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
